# New Reg Program Creation Step 1 Program Setting
# Update credential values on the "Credentials" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Credentials")

# "New Reg Program" url/username/password block (rows 2-4, column B)
$ws.Range("B2").Value = "http://testing.bsbtest.com/default.aspx?portalid=1141 "
$ws.Range("B3").Value = "bsbadmin1141"
$ws.Range("B4").Value = "Old4thWard#"

# "Login As Host" block (rows 7-9, column B)
$ws.Range("B7").Value = "https://stagingtshq.bsbtest.com/Default.aspx?portalid=24904"
$ws.Range("B8").Value = "athost"
$ws.Range("B9").Value = "ITW3546ctyz10@"

# Portal Range (rows 10-11, column B) - both bounds now the same single portal id
$ws.Range("B10").Value = "2152"
$ws.Range("B11").Value = "2152"

# Leave the last active selection on B11, matching the edited cell
$ws.Range("B11").Select()
